# Update "paises" (countries) COVID stats workbook to the newer snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 17:22"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1100608
$ws.Range("C4").Value = 5585
$ws.Range("E4").Value = 880497
$ws.Range("G4").Value = 166
$ws.Range("H4").Value = 64022

# --- Chile (row 29) ---
$ws.Range("B29").Value = 17008
$ws.Range("C29").Value = 985
$ws.Range("D29").Value = 9018
$ws.Range("E29").Value = 7756
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 234

# --- Republica Dominicana (row 47) ---
$ws.Range("B47").Value = 7288
$ws.Range("C47").Value = 316
$ws.Range("D47").Value = 1387
$ws.Range("E47").Value = 5588
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 313

# --- Moldavia (row 59) ---
$ws.Range("B59").Value = 3980
$ws.Range("C59").Value = 83
$ws.Range("E59").Value = 2586

# --- Republica de Chipre (row 94) ---
$ws.Range("B94").Value = 857
$ws.Range("C94").Value = 7
$ws.Range("E94").Value = 546

# --- Sri Lanka (row 103) ---
$ws.Range("B103").Value = 689
$ws.Range("C103").Value = 26
$ws.Range("E103").Value = 525

# --- Mali (row 112) ---
$ws.Range("B112").Value = 508
$ws.Range("C112").Value = 18
$ws.Range("D112").Value = 196
$ws.Range("E112").Value = 286

# --- Isla de Man / Guinea Ecuatorial swap places (row 127 / 128) ---
# Isla de Man overtakes Guinea Ecuatorial in total cases, so it now sits
# above it in the (descending) sorted list.
$ws.Range("A127").Value = "Isla de Man"
$ws.Range("B127").Value = 316
$ws.Range("C127").Value = 1
$ws.Range("D127").Value = 271
$ws.Range("E127").Value = 23
$ws.Range("F127").Value = 21
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 22

$ws.Range("A128").Value = "Guinea Ecuatorial"
$ws.Range("B128").Value = 315
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 9
$ws.Range("E128").Value = 305
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 1
